$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha), M (Volumen), N/O/P (Precio min/max/promedio), S (Precio $/Kg),
# derived from the updated weekly dataset. Keyed by row number.
$rowData = @{
    2 = @{ "D"=44748; "M"=300; "N"=2300; "O"=2300; "P"=2300; "S"=2300 }
    3 = @{ "D"=44749; "M"=120; "N"=2300; "O"=2300; "P"=2300; "S"=2300 }
    4 = @{ "D"=45075; "M"=240; "N"=3200; "O"=3200; "P"=3200; "S"=3200 }
    5 = @{ "D"=44763; "M"=50; "N"=2300; "O"=2300; "P"=2300; "S"=2300 }
    6 = @{ "D"=44762; "M"=50; "N"=2300; "O"=2300; "P"=2300; "S"=2300 }
    7 = @{ "D"=44753; "M"=160; "N"=2300; "O"=2300; "P"=2300; "S"=2300 }
    8 = @{ "D"=44811; "M"=60; "N"=2500; "O"=2500; "P"=2500; "S"=2500 }
    9 = @{ "D"=45062; "M"=60 }
    10 = @{ "D"=45068; "M"=50; "N"=3250; "O"=3250; "P"=3250; "S"=3250 }
    11 = @{ "D"=44760; "M"=80; "N"=2300; "O"=2300; "P"=2300; "S"=2300 }
    12 = @{ "D"=44435; "M"=130; "N"=1300; "O"=1300; "P"=1300; "S"=1300 }
    13 = @{ "D"=45054; "M"=25; "N"=2500; "O"=2500; "P"=2500; "S"=2500 }
    14 = @{ "D"=44417; "M"=80; "N"=1200; "O"=1200; "P"=1200; "S"=1200 }
    15 = @{ "D"=44476; "M"=80 }
    16 = @{ "D"=45076; "M"=100 }
    17 = @{ "D"=45044; "M"=150; "N"=3500; "O"=3500; "P"=3500; "S"=3500 }
    18 = @{ "D"=44812; "M"=50; "N"=2500; "O"=2500; "P"=2500; "S"=2500 }
    19 = @{ "D"=44357; "M"=35; "N"=1000; "O"=1000; "P"=1000; "S"=1000 }
    20 = @{ "D"=44438; "M"=60; "N"=1200; "O"=1200; "P"=1200; "S"=1200 }
    21 = @{ "D"=44343; "M"=60; "N"=1300; "O"=1300; "P"=1300; "S"=1300 }
    22 = @{ "D"=45079; "M"=30; "N"=2600; "O"=2600; "P"=2600; "S"=2600 }
    23 = @{ "D"=45042; "M"=25; "N"=3500; "O"=3500; "P"=3500; "S"=3500 }
    24 = @{ "D"=44431; "M"=100; "N"=1300; "O"=1300; "P"=1300; "S"=1300 }
    25 = @{ "D"=44418; "M"=40; "N"=1200; "O"=1200; "P"=1200; "S"=1200 }
    26 = @{ "D"=45041; "M"=80; "N"=3500; "O"=3500; "P"=3500; "S"=3500 }
    27 = @{ "D"=44405; "M"=50; "N"=1200; "O"=1200; "P"=1200; "S"=1200 }
    28 = @{ "D"=44432; "M"=30; "N"=1300; "O"=1300; "P"=1300; "S"=1300 }
    29 = @{ "D"=45055; "M"=25; "N"=2800; "O"=2800; "P"=2800; "S"=2800 }
    30 = @{ "D"=44424; "M"=50; "N"=1200; "O"=1200; "P"=1200; "S"=1200 }
    31 = @{ "D"=44473; "M"=120; "N"=1200; "O"=1200; "P"=1200; "S"=1200 }
    32 = @{ "D"=44830 }
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
